$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns updated per row: B, C, D, F, G, I, J, N
$cols = @("B","C","D","F","G","I","J","N")

$data = @(
    @(1.813677153944695,0.7065639160009596,0.1379950750549597,2.631178938091892,0.002501828940926806,1.094282161516929,0.3525354678427703,1.385982282758096),
    @(1.678304304843493,0.6523304722893499,0.135685494883262,2.591588512679365,0.002507634745315879,1.088857847933447,0.3416937574009893,1.406220596217732),
    @(1.596119825845676,0.619434472920716,0.1343179594307173,2.569109304773662,0.002511384081331946,1.086295189753933,0.3352687684444788,1.419257753582169),
    @(1.56286209936377,0.6061292223374721,0.1337734208037133,2.560406028647137,0.002512958541597698,1.085442839637238,0.3327084109077276,1.424723653769419),
    @(1.557353701558213,0.6039259048679924,0.1336837711561856,2.558988386332089,0.002513222797368886,1.085312868009893,0.3322867488126633,1.425640500049628),
    @(1.59567035925761,0.6192546299597552,0.1343105639535978,2.568990082016271,0.002511405126346039,1.086282918968564,0.3352340047786981,1.41933084909731),
    @(1.76680527682106,0.6877796364428264,0.1371882558100594,2.617146706372935,0.002503792580869421,1.092251788658281,0.3487488597494206,1.392833436538933),
    @(2.109929727720214,0.825427527885779,0.1432316732922914,2.726235936906448,0.002490320984899656,1.110103458579886,0.3771124302786149,1.345736472049065),
    @(2.366796260794786,0.9286578280064077,0.1479153204998056,2.815522570210447,0.002481300529408724,1.12704447132846,0.399117935598639,1.31412397256323),
    @(2.484730000745685,0.9761001246789078,0.1500988848146676,2.858171565121694,0.002477385028071699,1.135599354068844,0.4093894256744477,1.300397122701195),
    @(2.529547214163983,0.9941364126959229,0.1509333439109781,2.874617395170787,0.00247592917824956,1.138962234967636,0.4133170632429426,1.295293642437805),
    @(2.519887947411803,0.9902487959370205,0.1507532909742366,2.8710622910138,0.002476241529398566,1.138232471536924,0.4124694759160263,1.296388552726466),
    @(2.488413953959707,0.9775825505356011,0.1501673842677462,2.859518630040697,0.002477264716926643,1.135873540771655,0.4097117893306859,1.299975357544675),
    @(2.469155907623929,0.9698333923266773,0.1498094873529965,2.852486398892012,0.002477894942807059,1.134444729308512,0.4080275966934153,1.302184712941983),
    @(2.359111027907204,0.9255672006792679,0.147773682883809,2.812776535394192,0.002481560184598487,1.126502574599002,0.3984519660050552,1.315034258853416),
    @(2.291881353971064,0.8985359606805332,0.1465383273431655,2.788938661631761,0.00248385671515785,1.121848547441502,0.392644850539952,1.323084896097274),
    @(2.253314494212077,0.8830336087212345,0.1458327686977299,2.775418808871109,0.002485195319209573,1.119251441526529,0.389329294183355,1.327776917608187),
    @(2.240273842281056,0.8777925123298473,0.1455947356120788,2.77087393216695,0.002485651592338933,1.118385765594809,0.3882109038078028,1.329376102643547),
    @(2.299027507097662,0.9014087824817238,0.1466693172762632,2.791456445571129,0.00248361041479117,1.122335711287988,0.3932604841754994,1.322221522173315),
    @(2.497654313974977,0.9813009964764774,0.1503392734401103,2.862901233924958,0.002476963454142557,1.136563057902237,0.4105207520658922,1.298919254638854),
    @(2.628392073121972,1.033929301280807,0.1527820364080128,2.911319013969745,0.002472775799398533,1.146580965370859,0.4220233444967079,1.284241427256564),
    @(2.558529542754513,1.005802166964997,0.1514742477403814,2.885318613650242,0.002474996561349387,1.141167934831486,0.4158637111482903,1.292024605192971),
    @(2.295796468639708,0.9001098606673281,0.1466100821993308,2.79031757991396,0.002483721710267799,1.1221152199291,0.3929820845085601,1.322611655457564),
    @(2.016280723935665,0.7878290827597993,0.1415539703280047,2.69513455629999,0.002493810594230705,1.104608044901397,0.3692367972740129,1.35795422413592)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $r).Value = $rowVals[$j]
    }
}
